$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Insert 6 new rows before the current "Summe" row (row 34), which
# pushes it down to row 40 and frees up rows 34-39 for new content.
# Rows 32/33 already exist (currently just stub I32/I33 cells) and
# get overwritten below with real data.
# ---------------------------------------------------------------
$ws.Rows.Item(34).Resize(6).Insert()

# ---------------------------------------------------------------
# Row 32 - "Vorstellung"
# ---------------------------------------------------------------
$ws.Range("A32").Value = 12
$ws.Range("B32").Value = "Domaenenanalyse"
$ws.Range("C32").Value = "[SEMINAR]"
$ws.Range("D32").Value = "Themenfeldanalyse"
$ws.Range("E32").Value = "Vorstellung"

$ws.Range("F30").Copy()
$ws.Range("F32").PasteSpecial(-4122)
$ws.Range("G30").Copy()
$ws.Range("G32").PasteSpecial(-4122)
$ws.Range("H12").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("J30").Copy()
$ws.Range("J32").PasteSpecial(-4122)
$ws.Range("K30").Copy()
$ws.Range("K32").PasteSpecial(-4122)

$ws.Range("F32").Value = 44313
$ws.Range("G32").Value = 44317
$ws.Range("H32").Formula = "=ROUNDUP(((SUM(K32-J32)*24*60/60)/0.25),0)*0.25"
$ws.Range("J32").Value = 0.375
$ws.Range("K32").Value = 0.54166666666666663
$ws.Range("I32").Clear()

# ---------------------------------------------------------------
# Row 33 - "Notizen von der Vorstellung festhalten"
# ---------------------------------------------------------------
$ws.Range("A33").Value = 12
$ws.Range("B33").Value = "Domaenenanalyse"
$ws.Range("C33").Value = "[SEMINAR]"
$ws.Range("D33").Value = "Themenfeldanalyse"
$ws.Range("E33").Value = "Notizen von der Vorstellung festhalten"

$ws.Range("F30").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("G30").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("H12").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("J30").Copy()
$ws.Range("J33").PasteSpecial(-4122)
$ws.Range("K30").Copy()
$ws.Range("K33").PasteSpecial(-4122)

$ws.Range("F33").Value = 44313
$ws.Range("G33").Value = 44317
$ws.Range("H33").Formula = "=ROUNDUP(((SUM(K33-J33)*24*60/60)/0.25),0)*0.25"
$ws.Range("J33").Value = 0.58333333333333337
$ws.Range("K33").Value = 0.625
$ws.Range("I33").Clear()

# ---------------------------------------------------------------
# Row 34 - "Änderungen übernehmen und weitere Organisation"
# ---------------------------------------------------------------
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Domaenenanalyse"
$ws.Range("C34").Value = "[FEATURE]"
$ws.Range("D34").Value = "User Stories (mapping)"
$ws.Range("E34").Value = "Änderungen übernehmen und weitere Organisation"

$ws.Range("F30").Copy()
$ws.Range("F34").PasteSpecial(-4122)
$ws.Range("G30").Copy()
$ws.Range("G34").PasteSpecial(-4122)
$ws.Range("H12").Copy()
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("J30").Copy()
$ws.Range("J34").PasteSpecial(-4122)
$ws.Range("K30").Copy()
$ws.Range("K34").PasteSpecial(-4122)

$ws.Range("F34").Value = 44313
$ws.Range("G34").Value = 44317
$ws.Range("H34").ClearContents()
$ws.Range("I34").Formula = "=ROUNDUP(((SUM(K34-J34)*24*60/60)/0.25),0)*0.25"
$ws.Range("J34").Value = 0.70833333333333337
$ws.Range("K34").Value = 0.79166666666666663

# ---------------------------------------------------------------
# Rows 35-38: empty template rows (only date/time formatting kept)
# ---------------------------------------------------------------
$ws.Range("F30").Copy()
$ws.Range("F35:F38").PasteSpecial(-4122)
$ws.Range("G30").Copy()
$ws.Range("G35:G38").PasteSpecial(-4122)
$ws.Range("H12").Copy()
$ws.Range("H35:H38").PasteSpecial(-4122)
$ws.Range("J30").Copy()
$ws.Range("J35:J38").PasteSpecial(-4122)
$ws.Range("K30").Copy()
$ws.Range("K35:K38").PasteSpecial(-4122)

$ws.Range("F35:F38").ClearContents()
$ws.Range("G35:G38").ClearContents()
$ws.Range("H35:H38").ClearContents()
$ws.Range("J35:J38").ClearContents()
$ws.Range("K35:K38").ClearContents()
$ws.Range("I35:I38").Clear()

# Row 39 only keeps I39 (already carried along as an empty, styled
# cell by the row insert above).

# ---------------------------------------------------------------
# Data validation: extend prefix dropdown to the new rows.
# ---------------------------------------------------------------
$ws.Range("C34").Validation.Add(3, 1, 1, "=`$N`$3:`$N`$5")
$v1 = $ws.Range("C34").Validation
$v1.ErrorTitle = "Prefix nicht unterstützt"
$v1.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden`n"
$v1.InputTitle = "Prefix"
$v1.InputMessage = "Wählen Sie einen Prefix aus"

$ws.Range("C35:C39").Validation.Add(3, 1, 1, "=`$N`$3:`$N`$6")
$v2 = $ws.Range("C35:C39").Validation
$v2.ErrorTitle = "Prefix nicht unterstützt"
$v2.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden`n"
$v2.InputTitle = "Prefix"
$v2.InputMessage = "Wählen Sie einen Prefix aus"

# ---------------------------------------------------------------
# Sheet view: reflect the author having scrolled down to / selected A34.
# ---------------------------------------------------------------
$ws.Range("A34").Select()

Write-Host "edit applied"
